# Auto-generated script applying numeric cell updates per the commit diff.
# Updates currentAveragePrice / Leve Price / Leve Profit columns (H,I,J,K,L,M,N)
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 956.625
$ws.Range("I8").Value = 736.1429000000001
$ws.Range("K8").Value = 2208.4287
$ws.Range("M8").Value = -2069.4287
$ws.Range("H45").Value = 3000
$ws.Range("J45").Value = 3000
$ws.Range("L45").Value = 9000
$ws.Range("N45").Value = -9384
$ws.Range("H88").Value = 1860.5
$ws.Range("I88").Value = 1400
$ws.Range("J88").Value = 1911.6666
$ws.Range("K88").Value = 1400
$ws.Range("L88").Value = 1911.6666
$ws.Range("M88").Value = -994
$ws.Range("N88").Value = -2723.6666
$ws.Range("H91").Value = 1860.5
$ws.Range("I91").Value = 1400
$ws.Range("J91").Value = 1911.6666
$ws.Range("K91").Value = 1400
$ws.Range("L91").Value = 1911.6666
$ws.Range("M91").Value = 4
$ws.Range("N91").Value = -4719.6666
$ws.Range("H98").Value = 824.8182
$ws.Range("I98").Value = 788.09375
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 788.09375
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 709.90625
$ws.Range("N98").Value = -4996
$ws.Range("H112").Value = 1115.4889
$ws.Range("I112").Value = 1150
$ws.Range("J112").Value = 1113.8837
$ws.Range("K112").Value = 3450
$ws.Range("L112").Value = 3341.6511
$ws.Range("M112").Value = -2342
$ws.Range("N112").Value = -5557.6511
$ws.Range("H113").Value = 8966.087
$ws.Range("I113").Value = 3085.7144
$ws.Range("J113").Value = 11538.75
$ws.Range("K113").Value = 3085.7144
$ws.Range("L113").Value = 11538.75
$ws.Range("M113").Value = 168.2856000000002
$ws.Range("N113").Value = -18046.75
$ws.Range("H122").Value = 824.8182
$ws.Range("I122").Value = 788.09375
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2364.28125
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 85.71875
$ws.Range("N122").Value = -10900
$ws.Range("H129").Value = 1117
$ws.Range("I129").Value = 700
$ws.Range("K129").Value = 2100
$ws.Range("M129").Value = 2900
$ws.Range("H131").Value = 4353.56
$ws.Range("I131").Value = 254.6
$ws.Range("J131").Value = 4809
$ws.Range("K131").Value = 763.8
$ws.Range("L131").Value = 14427
$ws.Range("M131").Value = 4276.2
$ws.Range("N131").Value = -24507
$ws.Range("H132").Value = 1490.9672
$ws.Range("I132").Value = 1560.68
$ws.Range("J132").Value = 1174.091
$ws.Range("K132").Value = 4682.04
$ws.Range("L132").Value = 3522.273
$ws.Range("M132").Value = -2152.04
$ws.Range("N132").Value = -8582.272999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3465.5
$ws.Range("I122").Value = 2674.6667
$ws.Range("J122").Value = 3940
$ws.Range("K122").Value = 8024.000100000001
$ws.Range("L122").Value = 11820
$ws.Range("M122").Value = -5574.000100000001
$ws.Range("N122").Value = -16720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 30000
$ws.Range("J61").Value = 30000
$ws.Range("L61").Value = 30000
$ws.Range("N61").Value = -30626

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 35562.332
$ws.Range("I23").Value = 28454
$ws.Range("K23").Value = 28454
$ws.Range("M23").Value = -28214
$ws.Range("H27").Value = 35562.332
$ws.Range("I27").Value = 28454
$ws.Range("K27").Value = 28454
$ws.Range("M27").Value = -28262
$ws.Range("H31").Value = 2299.1633
$ws.Range("I31").Value = 1508.125
$ws.Range("J31").Value = 3788.1765
$ws.Range("K31").Value = 1508.125
$ws.Range("L31").Value = 3788.1765
$ws.Range("M31").Value = -1213.125
$ws.Range("N31").Value = -4378.1765
$ws.Range("H34").Value = 2299.1633
$ws.Range("I34").Value = 1508.125
$ws.Range("J34").Value = 3788.1765
$ws.Range("K34").Value = 1508.125
$ws.Range("L34").Value = 3788.1765
$ws.Range("M34").Value = -1306.125
$ws.Range("N34").Value = -4192.1765
$ws.Range("H58").Value = 1431.1228
$ws.Range("I58").Value = 984.7317
$ws.Range("J58").Value = 2575
$ws.Range("K58").Value = 984.7317
$ws.Range("L58").Value = 2575
$ws.Range("M58").Value = -781.7317
$ws.Range("N58").Value = -2981
$ws.Range("H99").Value = 8315.789000000001
$ws.Range("I99").Value = 2594.75
$ws.Range("J99").Value = 12476.546
$ws.Range("K99").Value = 2594.75
$ws.Range("L99").Value = 12476.546
$ws.Range("M99").Value = -1096.75
$ws.Range("N99").Value = -15472.546
$ws.Range("H126").Value = 8315.789000000001
$ws.Range("I126").Value = 2594.75
$ws.Range("J126").Value = 12476.546
$ws.Range("K126").Value = 7784.25
$ws.Range("L126").Value = 37429.638
$ws.Range("M126").Value = -5314.25
$ws.Range("N126").Value = -42369.638
$ws.Range("H134").Value = 1248.1506
$ws.Range("I134").Value = 1244.4464
$ws.Range("J134").Value = 1260.3529
$ws.Range("K134").Value = 3733.3392
$ws.Range("L134").Value = 3781.0587
$ws.Range("M134").Value = -1198.3392
$ws.Range("N134").Value = -8851.058700000001
$ws.Range("H136").Value = 1431.1228
$ws.Range("I136").Value = 984.7317
$ws.Range("J136").Value = 2575
$ws.Range("K136").Value = 2954.1951
$ws.Range("L136").Value = 7725
$ws.Range("M136").Value = -404.1950999999999
$ws.Range("N136").Value = -12825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2349.25
$ws.Range("I94").Value = 600
$ws.Range("J94").Value = 4098.5
$ws.Range("K94").Value = 1800
$ws.Range("L94").Value = 12295.5
$ws.Range("M94").Value = -1124
$ws.Range("N94").Value = -13647.5
$ws.Range("H130").Value = 1587.5
$ws.Range("I130").Value = 1500
$ws.Range("K130").Value = 4500
$ws.Range("M130").Value = 520
$ws.Range("H131").Value = 5630.968
$ws.Range("J131").Value = 6342.593
$ws.Range("L131").Value = 19027.779
$ws.Range("N131").Value = -29107.779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 3502.6667
$ws.Range("I6").Value = 508
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 508
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = -395
$ws.Range("N6").Value = -5226
$ws.Range("H16").Value = 3502.6667
$ws.Range("I16").Value = 508
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 508
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -258
$ws.Range("N16").Value = -5500
$ws.Range("H102").Value = 3400
$ws.Range("I102").Value = 1800
$ws.Range("K102").Value = 1800
$ws.Range("M102").Value = -178
$ws.Range("H122").Value = 2527.75
$ws.Range("I122").Value = 2017.7858
$ws.Range("J122").Value = 3037.7144
$ws.Range("K122").Value = 6053.357400000001
$ws.Range("L122").Value = 9113.143199999999
$ws.Range("M122").Value = -3603.357400000001
$ws.Range("N122").Value = -14013.1432
$ws.Range("H126").Value = 1929.2354
$ws.Range("I126").Value = 1720.3043
$ws.Range("J126").Value = 2366.0908
$ws.Range("K126").Value = 5160.9129
$ws.Range("L126").Value = 7098.2724
$ws.Range("M126").Value = -2690.9129
$ws.Range("N126").Value = -12038.2724
$ws.Range("H132").Value = 1762.6818
$ws.Range("I132").Value = 1409.6111
$ws.Range("J132").Value = 3351.5
$ws.Range("K132").Value = 4228.8333
$ws.Range("L132").Value = 10054.5
$ws.Range("M132").Value = -1698.8333
$ws.Range("N132").Value = -15114.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 170250
$ws.Range("I7").Value = 252875
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 252875
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -252763
$ws.Range("N7").Value = -5224
$ws.Range("H126").Value = 170250
$ws.Range("I126").Value = 252875
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 758625
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -756155
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 3000
$ws.Range("K14").Value = 3000
$ws.Range("M14").Value = -2832
$ws.Range("H126").Value = 79030.164
$ws.Range("I126").Value = 88677.69
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 266033.07
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -263563.07
$ws.Range("N126").Value = -10490
